# Update the "想去人数" (F column) figures on the "展览" and "全部类型"
# sheets to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F
$updates = @{
    2  = 6682
    5  = 52
    6  = 2090
    7  = 1591
    8  = 315
    9  = 1024
    10 = 477
    11 = 24
    12 = 5656
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
